$wb = $excel.ActiveWorkbook

# Remove the original (empty) Sheet1, keep Sheet2 which holds the combined dataset
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheet (formerly Sheet2) to Sheet1 so the dates/data line up
$wb.Worksheets.Item("Sheet2").Name = "Sheet1"
$wb.Worksheets.Item("Sheet1").Activate()
